# Auto-generated data-driven update of the cryptos worksheet.
# Parallel arrays (row number / column letter / new text value) are
# used instead of an array-of-arrays so each update is unambiguous.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowNums = @(2, 2, 3, 3, 4, 4, 5, 5, 6, 7, 7, 8, 8, 9, 9, 10, 10, 11, 11, 12, 12, 13, 13, 14, 14, 15, 15, 16, 16, 17, 17, 18, 18, 19, 19, 20, 20, 21, 21, 22, 22, 23, 23, 24, 24, 25, 25, 26, 26, 27, 27, 28, 28, 29, 29, 30, 30, 31, 31, 32, 32, 33, 33, 33, 33, 34, 34, 34, 34, 35, 35, 36, 36, 37, 37, 38, 38, 39, 39, 40, 40, 41, 42, 42, 43, 43, 44, 44, 45, 45, 46, 46, 47, 47, 48, 48, 48, 48, 49, 49, 49, 49, 50, 50, 51, 51)
$colLets = @('D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'B', 'C', 'D', 'E', 'B', 'C', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'D', 'E', 'B', 'C', 'D', 'E', 'B', 'C', 'D', 'E', 'D', 'E', 'D', 'E')
$newVals = @(
    '27.268.55'
    '  -4.50%  '
    '1.850.29'
    '  -6.02%  '
    '1.002'
    '  -0.93%  '
    '321.57'
    '  -0.31%  '
    '  -0.67%  '
    '0.4478'
    '  -6.62%  '
    '0.3830'
    '  -5.79%  '
    '47.94'
    '  -10.96%  '
    '0.07848'
    '  -8.02%  '
    '1.012'
    '  -4.90%  '
    '21.19'
    '  -5.89%  '
    '1.831.35'
    '  -7.01%  '
    '5.844'
    '  -5.96%  '
    '7.094'
    '  -7.06%  '
    '1.002'
    '  -0.97%  '
    '0.00001027'
    '  -4.79%  '
    '85.19'
    '  -6.69%  '
    '0.06522'
    '  -1.55%  '
    '16.90'
    '  -9.50%  '
    '1.001'
    '  -0.75%  '
    '5.468'
    '  -6.83%  '
    '27.248.21'
    '  -4.73%  '
    '10.75'
    '  -7.08%  '
    '2.252'
    '  -1.84%  '
    '2.044.20'
    '  -7.13%  '
    '150.89'
    '  -2.97%  '
    '19.39'
    '  -4.90%  '
    '5.479'
    '  -8.17%  '
    '2.030'
    '  -7.20%  '
    '119.46'
    '  -4.40%  '
    '0.09257'
    '  -4.12%  '
    'ImmutableX'
    'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    '0.9330'
    '  -6.09%  '
    'ARBITRUM'
    'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    '1.455'
    '  -0.95%  '
    '3.585'
    '  -2.83%  '
    '5.242'
    '  -7.71%  '
    '0.02220'
    '  -5.47%  '
    '0.05956'
    '  -5.13%  '
    '1.195'
    '  -5.05%  '
    '8.273'
    '  -9.87%  '
    '  -0.75%  '
    '0.5866'
    '  -6.17%  '
    '0.1859'
    '  -3.35%  '
    '10.05'
    '  -10.67%  '
    '1.259'
    '  -6.91%  '
    '0.5627'
    '  -6.10%  '
    '11.71'
    '  -10.71%  '
    'PancakeSwap'
    'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    '3.357'
    '  -1.68%  '
    'NEARProtocol'
    'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    '1.915'
    '  -7.97%  '
    '0.06840'
    '  -0.09%  '
    '108.44'
    '  -2.92%  '
)

if ($rowNums.Count -ne $colLets.Count -or $rowNums.Count -ne $newVals.Count) {
    throw "Mismatched update arrays: rows=$($rowNums.Count) cols=$($colLets.Count) vals=$($newVals.Count)"
}

for ($i = 0; $i -lt $rowNums.Count; $i++) {
    $row = $rowNums[$i]
    $col = $colLets[$i]
    $val = $newVals[$i]
    $cell = $ws.Range("$col$row")
    if ($col -eq 'D' -or $col -eq 'E') {
        # Columns D (Price) and E (Volume/1h) hold values that look numeric
        # (e.g. "27.268.55", "1.002", "0.00001027") or are percentage text
        # (e.g. "  -4.50%  "). The source workbook stores every one of
        # these as inline text (t="inlineStr"), never as a real number.
        # Forcing the cell to Text format before assigning the value stops
        # Excel's automatic type inference from turning these strings into
        # numbers/dates, then resetting the style back to Normal avoids
        # leaving a stray formatting difference behind.
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

